$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.8966636666666666
$ws.Cells.Item(2, 8).Value = 2.689991
$ws.Cells.Item(2, 9).Value = 0.334725143386341
$ws.Cells.Item(2, 10).Value = 0.3647360854412732
$ws.Cells.Item(2, 13).Value = 0.8151449999999999
$ws.Cells.Item(2, 14).Value = 2.445435
$ws.Cells.Item(2, 15).Value = 0.1271069095499719
$ws.Cells.Item(2, 16).Value = 0.1371035811308388
$ws.Cells.Item(2, 17).Value = 0.7309109045649999
$ws.Cells.Item(2, 18).Value = 6.578198141084999
$ws.Cells.Item(2, 19).Value = 0.04254587852450901
$ws.Cells.Item(2, 20).Value = 0.05000662348164217

$ws.Cells.Item(3, 7).Value = 0.8966636666666666
$ws.Cells.Item(3, 8).Value = 2.689991
$ws.Cells.Item(3, 9).Value = 0.334725143386341
$ws.Cells.Item(3, 10).Value = 0.3647360854412732
$ws.Cells.Item(3, 15).Value = 0.4802730342501803
$ws.Cells.Item(3, 16).Value = 0.5180454245123947
$ws.Cells.Item(3, 17).Value = 2.761744417709778
$ws.Cells.Item(3, 18).Value = 24.855699759388
$ws.Cells.Item(3, 19).Value = 0.1607594602539847
$ws.Cells.Item(3, 20).Value = 0.1889498602174134

$ws.Cells.Item(4, 7).Value = 0.8966636666666666
$ws.Cells.Item(4, 8).Value = 2.689991
$ws.Cells.Item(4, 9).Value = 0.334725143386341
$ws.Cells.Item(4, 10).Value = 0.3647360854412732
$ws.Cells.Item(4, 13).Value = 0.5185940000000001
$ws.Cells.Item(4, 14).Value = 1.555782
$ws.Cells.Item(4, 15).Value = 0.08086522109705406
$ws.Cells.Item(4, 16).Value = 0.08722508823947427
$ws.Cells.Item(4, 17).Value = 0.4650043975513334
$ws.Cells.Item(4, 18).Value = 4.185039577962001
$ws.Cells.Item(4, 19).Value = 0.02706762272667959
$ws.Cells.Item(4, 20).Value = 0.03181413723673548

$ws.Cells.Item(5, 7).Value = 0.8966636666666666
$ws.Cells.Item(5, 8).Value = 2.689991
$ws.Cells.Item(5, 9).Value = 0.334725143386341
$ws.Cells.Item(5, 10).Value = 0.3647360854412732
$ws.Cells.Item(5, 13).Value = 1.402793
$ws.Cells.Item(5, 14).Value = 2.805586
$ws.Cells.Item(5, 15).Value = 0.2187398352051889
$ws.Cells.Item(5, 16).Value = 0.1572954863942594
$ws.Cells.Item(5, 17).Value = 1.257833514954333
$ws.Cells.Item(5, 18).Value = 7.547001089726
$ws.Cells.Item(5, 19).Value = 0.07321772270336147
$ws.Cells.Item(5, 20).Value = 0.05737133996502321

$ws.Cells.Item(6, 7).Value = 0.8966636666666666
$ws.Cells.Item(6, 8).Value = 2.689991
$ws.Cells.Item(6, 9).Value = 0.334725143386341
$ws.Cells.Item(6, 10).Value = 0.3647360854412732
$ws.Cells.Item(6, 13).Value = 0.5965113333333333
$ws.Cells.Item(6, 14).Value = 1.789534
$ws.Cells.Item(6, 15).Value = 0.09301499989760488
$ws.Cells.Item(6, 16).Value = 0.1003304197230327
$ws.Cells.Item(6, 17).Value = 0.5348700393548889
$ws.Cells.Item(6, 18).Value = 4.813830354194
$ws.Cells.Item(6, 19).Value = 0.03113445917780629
$ws.Cells.Item(6, 20).Value = 0.03659412454045887

$ws.Cells.Item(7, 7).Value = 0.6827986666666667
$ws.Cells.Item(7, 9).Value = 0.2548891965854188
$ws.Cells.Item(7, 10).Value = 0.2777421703171357
$ws.Cells.Item(7, 13).Value = 0.8151449999999999
$ws.Cells.Item(7, 14).Value = 2.445435
$ws.Cells.Item(7, 15).Value = 0.1271069095499719
$ws.Cells.Item(7, 16).Value = 0.1371035811308388
$ws.Cells.Item(7, 17).Value = 0.5565799191399999
$ws.Cells.Item(7, 18).Value = 5.009219272259999
$ws.Cells.Item(7, 19).Value = 0.03239817805564783
$ws.Cells.Item(7, 20).Value = 0.03807944618153068

$ws.Cells.Item(8, 7).Value = 0.6827986666666667
$ws.Cells.Item(8, 9).Value = 0.2548891965854188
$ws.Cells.Item(8, 10).Value = 0.2777421703171357
$ws.Cells.Item(8, 15).Value = 0.4802730342501803
$ws.Cells.Item(8, 16).Value = 0.5180454245123947
$ws.Cells.Item(8, 18).Value = 18.927318330928
$ws.Cells.Item(8, 19).Value = 0.1224164078416698
$ws.Cells.Item(8, 20).Value = 0.1438830605269344

$ws.Cells.Item(9, 7).Value = 0.6827986666666667
$ws.Cells.Item(9, 9).Value = 0.2548891965854188
$ws.Cells.Item(9, 10).Value = 0.2777421703171357
$ws.Cells.Item(9, 13).Value = 0.5185940000000001
$ws.Cells.Item(9, 14).Value = 1.555782
$ws.Cells.Item(9, 15).Value = 0.08086522109705406
$ws.Cells.Item(9, 16).Value = 0.08722508823947427
$ws.Cells.Item(9, 17).Value = 0.3540952917413334
$ws.Cells.Item(9, 18).Value = 3.186857625672
$ws.Cells.Item(9, 19).Value = 0.02061167123713037
$ws.Cells.Item(9, 20).Value = 0.02422608531373525

$ws.Cells.Item(10, 7).Value = 0.6827986666666667
$ws.Cells.Item(10, 9).Value = 0.2548891965854188
$ws.Cells.Item(10, 10).Value = 0.2777421703171357
$ws.Cells.Item(10, 13).Value = 1.402793
$ws.Cells.Item(10, 14).Value = 2.805586
$ws.Cells.Item(10, 15).Value = 0.2187398352051889
$ws.Cells.Item(10, 16).Value = 0.1572954863942594
$ws.Cells.Item(10, 17).Value = 0.9578251900093333
$ws.Cells.Item(10, 18).Value = 5.746951140056
$ws.Cells.Item(10, 19).Value = 0.05575442085667751
$ws.Cells.Item(10, 20).Value = 0.04368758977223108

$ws.Cells.Item(11, 7).Value = 0.6827986666666667
$ws.Cells.Item(11, 9).Value = 0.2548891965854188
$ws.Cells.Item(11, 10).Value = 0.2777421703171357
$ws.Cells.Item(11, 13).Value = 0.5965113333333333
$ws.Cells.Item(11, 14).Value = 1.789534
$ws.Cells.Item(11, 15).Value = 0.09301499989760488
$ws.Cells.Item(11, 16).Value = 0.1003304197230327
$ws.Cells.Item(11, 17).Value = 0.4072971430515555
$ws.Cells.Item(11, 18).Value = 3.665674287464
$ws.Cells.Item(11, 19).Value = 0.02370851859429332
$ws.Cells.Item(11, 20).Value = 0.02786598852270427

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.221369
$ws.Cells.Item(12, 8).Value = 0.664107
$ws.Cells.Item(12, 9).Value = 0.08263719499391366
$ws.Cells.Item(12, 10).Value = 0.09004631892602898
$ws.Cells.Item(12, 13).Value = 0.8151449999999999
$ws.Cells.Item(12, 14).Value = 2.445435
$ws.Cells.Item(12, 15).Value = 0.1271069095499719
$ws.Cells.Item(12, 16).Value = 0.1371035811308388
$ws.Cells.Item(12, 17).Value = 0.180447833505
$ws.Cells.Item(12, 18).Value = 1.624030501545
$ws.Cells.Item(12, 19).Value = 0.01050375846955477
$ws.Cells.Item(12, 20).Value = 0.0123456727924082

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.221369
$ws.Cells.Item(13, 8).Value = 0.664107
$ws.Cells.Item(13, 9).Value = 0.08263719499391366
$ws.Cells.Item(13, 10).Value = 0.09004631892602898
$ws.Cells.Item(13, 15).Value = 0.4802730342501803
$ws.Cells.Item(13, 16).Value = 0.5180454245123947
$ws.Cells.Item(13, 17).Value = 0.6818215376973334
$ws.Cells.Item(13, 18).Value = 6.136393839276001
$ws.Cells.Item(13, 19).Value = 0.03968841638165072
$ws.Cells.Item(13, 20).Value = 0.04664808351381316

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.221369
$ws.Cells.Item(14, 8).Value = 0.664107
$ws.Cells.Item(14, 9).Value = 0.08263719499391366
$ws.Cells.Item(14, 10).Value = 0.09004631892602898
$ws.Cells.Item(14, 13).Value = 0.5185940000000001
$ws.Cells.Item(14, 14).Value = 1.555782
$ws.Cells.Item(14, 15).Value = 0.08086522109705406
$ws.Cells.Item(14, 16).Value = 0.08722508823947427
$ws.Cells.Item(14, 17).Value = 0.114800635186
$ws.Cells.Item(14, 18).Value = 1.033205716674
$ws.Cells.Item(14, 19).Value = 0.006682475044023197
$ws.Cells.Item(14, 20).Value = 0.007854298113962719

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.221369
$ws.Cells.Item(15, 8).Value = 0.664107
$ws.Cells.Item(15, 9).Value = 0.08263719499391366
$ws.Cells.Item(15, 10).Value = 0.09004631892602898
$ws.Cells.Item(15, 13).Value = 1.402793
$ws.Cells.Item(15, 14).Value = 2.805586
$ws.Cells.Item(15, 15).Value = 0.2187398352051889
$ws.Cells.Item(15, 16).Value = 0.1572954863942594
$ws.Cells.Item(15, 17).Value = 0.310534883617
$ws.Cells.Item(15, 18).Value = 1.863209301702
$ws.Cells.Item(15, 19).Value = 0.01807604641478774
$ws.Cells.Item(15, 20).Value = 0.01416387953348233

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.221369
$ws.Cells.Item(16, 8).Value = 0.664107
$ws.Cells.Item(16, 9).Value = 0.08263719499391366
$ws.Cells.Item(16, 10).Value = 0.09004631892602898
$ws.Cells.Item(16, 13).Value = 0.5965113333333333
$ws.Cells.Item(16, 14).Value = 1.789534
$ws.Cells.Item(16, 15).Value = 0.09301499989760488
$ws.Cells.Item(16, 16).Value = 0.1003304197230327
$ws.Cells.Item(16, 17).Value = 0.1320491173486667
$ws.Cells.Item(16, 18).Value = 1.188442056138
$ws.Cells.Item(16, 19).Value = 0.007686498683897233
$ws.Cells.Item(16, 20).Value = 0.009034384972362553

$ws.Cells.Item(17, 7).Value = 0.6612465000000001
$ws.Cells.Item(17, 8).Value = 1.322493
$ws.Cells.Item(17, 9).Value = 0.2468437584284291
$ws.Cells.Item(17, 10).Value = 0.1793169270244717
$ws.Cells.Item(17, 13).Value = 0.8151449999999999
$ws.Cells.Item(17, 14).Value = 2.445435
$ws.Cells.Item(17, 15).Value = 0.1271069095499719
$ws.Cells.Item(17, 16).Value = 0.1371035811308388
$ws.Cells.Item(17, 17).Value = 0.5390117782425
$ws.Cells.Item(17, 18).Value = 3.234070669455
$ws.Cells.Item(17, 19).Value = 0.03137554727553744
$ws.Cells.Item(17, 20).Value = 0.02458499285243237

$ws.Cells.Item(18, 7).Value = 0.6612465000000001
$ws.Cells.Item(18, 8).Value = 1.322493
$ws.Cells.Item(18, 9).Value = 0.2468437584284291
$ws.Cells.Item(18, 10).Value = 0.1793169270244717
$ws.Cells.Item(18, 15).Value = 0.4802730342501803
$ws.Cells.Item(18, 16).Value = 0.5180454245123947
$ws.Cells.Item(18, 17).Value = 2.036654208254
$ws.Cells.Item(18, 18).Value = 12.219925249524
$ws.Cells.Item(18, 19).Value = 0.1185524008461402
$ws.Cells.Item(18, 20).Value = 0.09289431358265056

$ws.Cells.Item(19, 7).Value = 0.6612465000000001
$ws.Cells.Item(19, 8).Value = 1.322493
$ws.Cells.Item(19, 9).Value = 0.2468437584284291
$ws.Cells.Item(19, 10).Value = 0.1793169270244717
$ws.Cells.Item(19, 13).Value = 0.5185940000000001
$ws.Cells.Item(19, 14).Value = 1.555782
$ws.Cells.Item(19, 15).Value = 0.08086522109705406
$ws.Cells.Item(19, 16).Value = 0.08722508823947427
$ws.Cells.Item(19, 17).Value = 0.3429184674210001
$ws.Cells.Item(19, 18).Value = 2.057510804526001
$ws.Cells.Item(19, 19).Value = 0.01996107510174272
$ws.Cells.Item(19, 20).Value = 0.01564093478254092

$ws.Cells.Item(20, 7).Value = 0.6612465000000001
$ws.Cells.Item(20, 8).Value = 1.322493
$ws.Cells.Item(20, 9).Value = 0.2468437584284291
$ws.Cells.Item(20, 10).Value = 0.1793169270244717
$ws.Cells.Item(20, 13).Value = 1.402793
$ws.Cells.Item(20, 14).Value = 2.805586
$ws.Cells.Item(20, 15).Value = 0.2187398352051889
$ws.Cells.Item(20, 16).Value = 0.1572954863942594
$ws.Cells.Item(20, 17).Value = 0.9275919614745001
$ws.Cells.Item(20, 18).Value = 3.710367845898
$ws.Cells.Item(20, 19).Value = 0.05399456304006405
$ws.Cells.Item(20, 20).Value = 0.02820574325503819

$ws.Cells.Item(21, 7).Value = 0.6612465000000001
$ws.Cells.Item(21, 8).Value = 1.322493
$ws.Cells.Item(21, 9).Value = 0.2468437584284291
$ws.Cells.Item(21, 10).Value = 0.1793169270244717
$ws.Cells.Item(21, 13).Value = 0.5965113333333333
$ws.Cells.Item(21, 14).Value = 1.789534
$ws.Cells.Item(21, 15).Value = 0.09301499989760488
$ws.Cells.Item(21, 16).Value = 0.1003304197230327
$ws.Cells.Item(21, 17).Value = 0.394441031377
$ws.Cells.Item(21, 18).Value = 2.366646188262
$ws.Cells.Item(21, 19).Value = 0.02296017216494474
$ws.Cells.Item(21, 20).Value = 0.01799094255180968

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 0.216728
$ws.Cells.Item(22, 8).Value = 0.650184
$ws.Cells.Item(22, 9).Value = 0.08090470660589748
$ws.Cells.Item(22, 10).Value = 0.08815849829109049
$ws.Cells.Item(22, 13).Value = 0.8151449999999999
$ws.Cells.Item(22, 14).Value = 2.445435
$ws.Cells.Item(22, 15).Value = 0.1271069095499719
$ws.Cells.Item(22, 16).Value = 0.1371035811308388
$ws.Cells.Item(22, 17).Value = 0.17666474556
$ws.Cells.Item(22, 18).Value = 1.58998271004
$ws.Cells.Item(22, 19).Value = 0.01028354722472282
$ws.Cells.Item(22, 20).Value = 0.01208684582282544

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 0.216728
$ws.Cells.Item(23, 8).Value = 0.650184
$ws.Cells.Item(23, 9).Value = 0.08090470660589748
$ws.Cells.Item(23, 10).Value = 0.08815849829109049
$ws.Cells.Item(23, 15).Value = 0.4802730342501803
$ws.Cells.Item(23, 16).Value = 0.5180454245123947
$ws.Cells.Item(23, 17).Value = 0.6675271525013333
$ws.Cells.Item(23, 18).Value = 6.007744372512001
$ws.Cells.Item(23, 19).Value = 0.03885634892673499
$ws.Cells.Item(23, 20).Value = 0.04567010667158319

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 0.216728
$ws.Cells.Item(24, 8).Value = 0.650184
$ws.Cells.Item(24, 9).Value = 0.08090470660589748
$ws.Cells.Item(24, 10).Value = 0.08815849829109049
$ws.Cells.Item(24, 13).Value = 0.5185940000000001
$ws.Cells.Item(24, 14).Value = 1.555782
$ws.Cells.Item(24, 15).Value = 0.08086522109705406
$ws.Cells.Item(24, 16).Value = 0.08722508823947427
$ws.Cells.Item(24, 17).Value = 0.112393840432
$ws.Cells.Item(24, 18).Value = 1.011544563888
$ws.Cells.Item(24, 19).Value = 0.006542376987478189
$ws.Cells.Item(24, 20).Value = 0.007689632792499909

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 0.216728
$ws.Cells.Item(25, 8).Value = 0.650184
$ws.Cells.Item(25, 9).Value = 0.08090470660589748
$ws.Cells.Item(25, 10).Value = 0.08815849829109049
$ws.Cells.Item(25, 13).Value = 1.402793
$ws.Cells.Item(25, 14).Value = 2.805586
$ws.Cells.Item(25, 15).Value = 0.2187398352051889
$ws.Cells.Item(25, 16).Value = 0.1572954863942594
$ws.Cells.Item(25, 17).Value = 0.304024521304
$ws.Cells.Item(25, 18).Value = 1.824147127824
$ws.Cells.Item(25, 19).Value = 0.01769708219029817
$ws.Cells.Item(25, 20).Value = 0.01386693386848456

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 0.216728
$ws.Cells.Item(26, 8).Value = 0.650184
$ws.Cells.Item(26, 9).Value = 0.08090470660589748
$ws.Cells.Item(26, 10).Value = 0.08815849829109049
$ws.Cells.Item(26, 13).Value = 0.5965113333333333
$ws.Cells.Item(26, 14).Value = 1.789534
$ws.Cells.Item(26, 15).Value = 0.09301499989760488
$ws.Cells.Item(26, 16).Value = 0.1003304197230327
$ws.Cells.Item(26, 17).Value = 0.1292807082506666
$ws.Cells.Item(26, 18).Value = 1.163526374256
$ws.Cells.Item(26, 19).Value = 0.007525351276663306
$ws.Cells.Item(26, 20).Value = 0.008844979135697373
